$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 572 metrics (K,L,M,N) for the existing mortality/DSVA_DPE run ---
$ws.Range("K572").Value = 63
$ws.Range("L572").Value = 0.493
$ws.Range("M572").Value = 0.5570000000000001
$ws.Range("N572").Value = 0.16

# --- Append new experiment run rows 573-582 ---
# Row 573
$ws.Range("A573").Value = 'outputs/2024-04-29/06-40-43'
$ws.Range("B573").Value = $false
$ws.Range("C573").Value = 'mimiciii'
$ws.Range("D573").Value = 'mortality'
$ws.Range("E573").Value = 'VC'
$ws.Range("F573").Value = 'descemb_bert'
$ws.Range("G573").Value = 'ehr_model'
$ws.Range("H573").Value = "'False"
$ws.Range("H573").ClearFormats()
$ws.Range("I573").Value = "'True"
$ws.Range("I573").ClearFormats()
$ws.Range("J573").Value = "'False"
$ws.Range("J573").ClearFormats()
$ws.Range("K573").Value = 1
$ws.Range("K573").Value = ""
$ws.Range("K573").ClearFormats()
$ws.Range("L573").Value = 1
$ws.Range("L573").Value = ""
$ws.Range("L573").ClearFormats()
$ws.Range("M573").Value = 1
$ws.Range("M573").Value = ""
$ws.Range("M573").ClearFormats()
$ws.Range("N573").Value = 1
$ws.Range("N573").Value = ""
$ws.Range("N573").ClearFormats()

# Row 574
$ws.Range("A574").Value = 'outputs/2024-04-29/06-45-17'
$ws.Range("B574").Value = $true
$ws.Range("C574").Value = 'eicu'
$ws.Range("D574").Value = 'mlm'
$ws.Range("E574").Value = 'NV'
$ws.Range("F574").Value = 'None'
$ws.Range("G574").Value = 'descemb_bert'
$ws.Range("H574").Value = "'True"
$ws.Range("H574").ClearFormats()
$ws.Range("I574").Value = "'False"
$ws.Range("I574").ClearFormats()
$ws.Range("J574").Value = "'False"
$ws.Range("J574").ClearFormats()
$ws.Range("K574").Value = 1000
$ws.Range("L574").Value = 0.555
$ws.Range("M574").Value = 0
$ws.Range("N574").Value = 0

# Row 575
$ws.Range("A575").Value = 'outputs/2024-04-29/07-26-37'
$ws.Range("B575").Value = $true
$ws.Range("C575").Value = 'mimiciii'
$ws.Range("D575").Value = 'diagnosis'
$ws.Range("E575").Value = 'VA'
$ws.Range("F575").Value = 'descemb_bert'
$ws.Range("G575").Value = 'ehr_model'
$ws.Range("H575").Value = "'False"
$ws.Range("H575").ClearFormats()
$ws.Range("I575").Value = "'False"
$ws.Range("I575").ClearFormats()
$ws.Range("J575").Value = "'False"
$ws.Range("J575").ClearFormats()
$ws.Range("K575").Value = 223
$ws.Range("L575").Value = 1.128
$ws.Range("M575").Value = 0.767
$ws.Range("N575").Value = 0.637

# Row 576
$ws.Range("A576").Value = 'outputs/2024-04-29/08-43-07'
$ws.Range("B576").Value = $true
$ws.Range("C576").Value = 'mimiciii'
$ws.Range("D576").Value = 'diagnosis'
$ws.Range("E576").Value = 'DSVA'
$ws.Range("F576").Value = 'descemb_bert'
$ws.Range("G576").Value = 'ehr_model'
$ws.Range("H576").Value = "'False"
$ws.Range("H576").ClearFormats()
$ws.Range("I576").Value = "'False"
$ws.Range("I576").ClearFormats()
$ws.Range("J576").Value = "'False"
$ws.Range("J576").ClearFormats()
$ws.Range("K576").Value = 333
$ws.Range("L576").Value = 1.141
$ws.Range("M576").Value = 0.768
$ws.Range("N576").Value = 0.638

# Row 577
$ws.Range("A577").Value = 'outputs/2024-04-29/10-47-28'
$ws.Range("B577").Value = $false
$ws.Range("C577").Value = 'mimiciii'
$ws.Range("D577").Value = 'diagnosis'
$ws.Range("E577").Value = 'DSVA_DPE'
$ws.Range("F577").Value = 'descemb_bert'
$ws.Range("G577").Value = 'ehr_model'
$ws.Range("H577").Value = "'False"
$ws.Range("H577").ClearFormats()
$ws.Range("I577").Value = "'False"
$ws.Range("I577").ClearFormats()
$ws.Range("J577").Value = "'False"
$ws.Range("J577").ClearFormats()
$ws.Range("K577").Value = 1
$ws.Range("K577").Value = ""
$ws.Range("K577").ClearFormats()
$ws.Range("L577").Value = 1
$ws.Range("L577").Value = ""
$ws.Range("L577").ClearFormats()
$ws.Range("M577").Value = 1
$ws.Range("M577").Value = ""
$ws.Range("M577").ClearFormats()
$ws.Range("N577").Value = 1
$ws.Range("N577").Value = ""
$ws.Range("N577").ClearFormats()

# Row 578
$ws.Range("A578").Value = 'outputs/2024-04-29/10-47-31'
$ws.Range("B578").Value = $true
$ws.Range("C578").Value = 'mimiciii'
$ws.Range("D578").Value = 'diagnosis'
$ws.Range("E578").Value = 'VC'
$ws.Range("F578").Value = 'descemb_bert'
$ws.Range("G578").Value = 'ehr_model'
$ws.Range("H578").Value = "'False"
$ws.Range("H578").ClearFormats()
$ws.Range("I578").Value = "'False"
$ws.Range("I578").ClearFormats()
$ws.Range("J578").Value = "'False"
$ws.Range("J578").ClearFormats()
$ws.Range("K578").Value = 269
$ws.Range("L578").Value = 1.104
$ws.Range("M578").Value = 0.768
$ws.Range("N578").Value = 0.638

# Row 579
$ws.Range("A579").Value = 'outputs/2024-04-29/12-09-48'
$ws.Range("B579").Value = $true
$ws.Range("C579").Value = 'mimiciii'
$ws.Range("D579").Value = 'los_3day'
$ws.Range("E579").Value = 'VA'
$ws.Range("F579").Value = 'descemb_bert'
$ws.Range("G579").Value = 'ehr_model'
$ws.Range("H579").Value = "'False"
$ws.Range("H579").ClearFormats()
$ws.Range("I579").Value = "'False"
$ws.Range("I579").ClearFormats()
$ws.Range("J579").Value = "'False"
$ws.Range("J579").ClearFormats()
$ws.Range("K579").Value = 195
$ws.Range("L579").Value = 1.683
$ws.Range("M579").Value = 0.512
$ws.Range("N579").Value = 0.348

# Row 580
$ws.Range("A580").Value = 'outputs/2024-04-29/13-15-34'
$ws.Range("B580").Value = $true
$ws.Range("C580").Value = 'mimiciii'
$ws.Range("D580").Value = 'los_3day'
$ws.Range("E580").Value = 'DSVA'
$ws.Range("F580").Value = 'descemb_bert'
$ws.Range("G580").Value = 'ehr_model'
$ws.Range("H580").Value = "'False"
$ws.Range("H580").ClearFormats()
$ws.Range("I580").Value = "'False"
$ws.Range("I580").ClearFormats()
$ws.Range("J580").Value = "'False"
$ws.Range("J580").ClearFormats()
$ws.Range("K580").Value = 191
$ws.Range("L580").Value = 1.669
$ws.Range("M580").Value = 0.506
$ws.Range("N580").Value = 0.35

# Row 581
$ws.Range("A581").Value = 'outputs/2024-04-29/14-26-06'
$ws.Range("B581").Value = $false
$ws.Range("C581").Value = 'mimiciii'
$ws.Range("D581").Value = 'los_3day'
$ws.Range("E581").Value = 'DSVA_DPE'
$ws.Range("F581").Value = 'descemb_bert'
$ws.Range("G581").Value = 'ehr_model'
$ws.Range("H581").Value = "'False"
$ws.Range("H581").ClearFormats()
$ws.Range("I581").Value = "'False"
$ws.Range("I581").ClearFormats()
$ws.Range("J581").Value = "'False"
$ws.Range("J581").ClearFormats()
$ws.Range("K581").Value = 1
$ws.Range("K581").Value = ""
$ws.Range("K581").ClearFormats()
$ws.Range("L581").Value = 1
$ws.Range("L581").Value = ""
$ws.Range("L581").ClearFormats()
$ws.Range("M581").Value = 1
$ws.Range("M581").Value = ""
$ws.Range("M581").ClearFormats()
$ws.Range("N581").Value = 1
$ws.Range("N581").Value = ""
$ws.Range("N581").ClearFormats()

# Row 582
$ws.Range("A582").Value = 'outputs/2024-04-29/14-26-09'
$ws.Range("B582").Value = $false
$ws.Range("C582").Value = 'mimiciii'
$ws.Range("D582").Value = 'los_3day'
$ws.Range("E582").Value = 'VC'
$ws.Range("F582").Value = 'descemb_bert'
$ws.Range("G582").Value = 'ehr_model'
$ws.Range("H582").Value = "'False"
$ws.Range("H582").ClearFormats()
$ws.Range("I582").Value = "'False"
$ws.Range("I582").ClearFormats()
$ws.Range("J582").Value = "'False"
$ws.Range("J582").ClearFormats()
$ws.Range("K582").Value = 61
$ws.Range("L582").Value = 1.525
$ws.Range("M582").Value = 0.511
$ws.Range("N582").Value = 0.35

